$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: Serveur -> Client, BS, battleState, description
$ws.Range("A24").Value = "Serveur"
$ws.Range("B24").Value = "Client"
$ws.Range("C24").Value = "BS"
$ws.Range("D24").Value = "battleState"
$ws.Range("E24").Value = "Indique au client que l'état du combat a changé (phase préparation, phase combat, etc…)"

# Row 25: Serveur -> Client, Cs, idPerso;[1 ou 0], description
$ws.Range("A25").Value = "Serveur"
$ws.Range("B25").Value = "Client"
$ws.Range("C25").Value = "Cs"
$ws.Range("E25").Value = "Indique au client l'état prêt/pas prêt d'un personnage"
$ws.Range("D25").Value = "idPerso;[1 ou 0]"

# Row 26: Client -> Serveur, Cs, (no params), description
$ws.Range("A26").Value = "Client"
$ws.Range("B26").Value = "Serveur"
$ws.Range("C26").Value = "Cs"
$ws.Range("E26").Value = "Indique au serveur que la position est verrouillée (joueur prêt à commencer le combat)"

# Row 27: Client -> Serveur, CP, cellX;cellY, description
$ws.Range("A27").Value = "Client"
$ws.Range("B27").Value = "Serveur"
$ws.Range("C27").Value = "CP"
$ws.Range("D27").Value = "cellX;cellY"
$ws.Range("E27").Value = "Demande un changement de position de départ au serveur"

# Row 28: Serveur -> Client, CP, idPerso;cellX;cellY, description
$ws.Range("A28").Value = "Serveur"
$ws.Range("B28").Value = "Client"
$ws.Range("C28").Value = "CP"
$ws.Range("D28").Value = "idPerso;cellX;cellY"
$ws.Range("E28").Value = "Indique un changement de position du personnage (position de départ ou téléportation en combat)"

# Row heights
$ws.Rows.Item(24).RowHeight = 30
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 45

# Update sheet view: pane/selection to match new state
$ws.Application.ActiveWindow.ScrollRow = 23
$ws.Range("A29").Select()
